$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 6183
$ws1.Range("F12").Value = 1252
$ws1.Range("F13").Value = 1252
$ws1.Range("F22").Value = 4509
$ws1.Range("F26").Value = 55

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 6183
$ws4.Range("F12").Value = 1252
$ws4.Range("F13").Value = 1252
$ws4.Range("F22").Value = 4509
$ws4.Range("F27").Value = 55
